$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new run ("Run 50") was added to the results table. The former last data
# column (AZ), which held the "Mean" header/values, is pushed one column to
# the right (BA); the new "Run 50" data takes over column AZ.

# 1) Create the new "Mean" column in BA, copying AZ1's header formatting.
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BA1").Value = "Mean"

# 2) Move the mean values (old AZ data) into BA, using the new recalculated mean.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 53).Value = 0.67291612
}

# 3) Turn the old "Mean" column (AZ) into the new "Run 50" column.
$ws.Range("AZ1").Value = "Run 50"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = 0.35741612
}
